$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "nm_ibu" column (K) entirely — columns to the right shift left.
$ws.Columns("K").Delete() | Out-Null

# Remove the "tp_id" column (now at L after the previous delete) entirely.
$ws.Columns("L").Delete() | Out-Null

# Fix the sample tgl_lhr value typo; force text so Excel doesn't coerce it to a date.
$ws.Range("E2").Value = "'2004-02-10"

# Replace the old tp_id note with the new tgl_lahir formatting note.
$ws.Range("L3").Value = "Ket: tgl_lahir pengisian= thn-bln-tgl, contoh: '2021-01-15"

# Move the selection to where the author left off.
$ws.Range("I10").Select() | Out-Null
